# Commit: "added temp check inputs to gui"
#
# On the "summary" sheet, split the old "Surface Temperature" and
# "Air Temperature" header columns into a dedicated *Check column
# (new input) plus the existing Min/Max columns, and make "summary"
# the active tab/selection instead of "Stations & Drops".

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("summary")

# Rename the header cells that used to hold the bare "Surface
# Temperature" / "Air Temperature" labels into the new "* Check"
# labels (the Min/Max columns next to them already exist and are
# untouched).
$summary.Range("L1").Value = "Surface Temp Check"
$summary.Range("O1").Value = "Air Temp Check"

# Make "summary" the active sheet/tab (it was "Stations & Drops"
# before), and move the selection on summary from T2 to O2.
$summary.Activate()
$summary.Range("O2").Select()
